# ringihistories schema fix:
#  - id row's ColumnName comment was wrong (duplicate of the next row) -> clear it
#  - ringiNo -> ringino, and its type corrected from varchar(255) to int(10)
#  - ringiseq type corrected from tinyint(1) to int(3)
#  - approverLayer -> approverlayer, type corrected from tinyint(1) to int(3)
#  - ringiAction -> ringiaction
# Order below mirrors the order the cells were touched in the original
# edit session so the shared-string table comes out in the same order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").ClearContents()

$ws.Range("C3").Value = "int(10)"
$ws.Range("C4").Value = "int(3)"

$ws.Range("A8").Value = "ringiaction"
$ws.Range("A7").Value = "approverlayer"
$ws.Range("C7").Value = "int(3)"

$ws.Range("A3").Value = "ringino"

# Column widths for the first three columns (set on the final, wider
# content so the widths are meaningful).
$ws.Columns.Item(1).ColumnWidth = 24.33
$ws.Columns.Item(2).ColumnWidth = 16.83
$ws.Columns.Item(3).ColumnWidth = 19

# Active cell moved from D3 to C3.
$ws.Range("C3").Select()
